# Consolidate "Slide" + " " runs into a single "Slide " run, and
# "an" + " " runs into a single "an " run, on both slides — matching
# the PowerPoint writer behaviour of merging adjacent same-formatted
# text runs.

$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Characters(1, 6).Text = "Slide "
$s1.Shapes.Item(3).TextFrame.TextRange.Characters(1, 3).Text = "an "

$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Characters(1, 6).Text = "Slide "
$s2.Shapes.Item(4).TextFrame.TextRange.Characters(1, 3).Text = "an "
